$d = $word.ActiveDocument
Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
Write-Output ("Bookmarks: " + $d.Bookmarks.Count)
foreach ($bm in $d.Bookmarks) {
    Write-Output ("Bookmark: " + $bm.Name)
}
